$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '59.198.61'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.522.13'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').Value = '  +5.60%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '135.05'
$ws.Range('E6').Value = '  +4.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.567'
$ws.Range('E8').Value = '  +3.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.520.32'
$ws.Range('E9').Value = '  +2.60%  '
$ws.Range('E10').Value = '  +4.57%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.25'
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.334'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.967.40'
$ws.Range('E14').Value = '  +3.23%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '59.148.75'
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '22.48'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('E17').Value = '  +3.40%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.518.14'
$ws.Range('E18').Value = '  +2.78%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.70'
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '321.93'
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.19'
$ws.Range('E22').Value = '  +9.23%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  +3.82%  '
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.53'
$ws.Range('E28').Value = '  +3.58%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0768'
$ws.Range('E29').Value = '  +6.14%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '173.28'
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('E31').Value = '  +5.47%  '
$ws.Range('E32').Value = '  +4.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.32'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('E35').Value = '  +0.27%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '18.18'
$ws.Range('E36').Value = '  +2.59%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  +0.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.53'
$ws.Range('E39').Value = '  +4.77%  '
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.798'
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.50'
$ws.Range('E42').Value = '  +3.39%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '278.64'
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('E44').Value = '  +3.32%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '132.17'
$ws.Range('E45').Value = '  +10.00%  '
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0935'
$ws.Range('E47').Value = '  +3.03%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0513'
$ws.Range('E48').Value = '  +5.75%  '
$ws.Range('E49').Value = '  +5.06%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '17.17'
$ws.Range('E50').Value = '  +2.97%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.760.77'
$ws.Range('E51').Value = '  +3.19%  '
